$d = $word.ActiveDocument

# Update the certificate date
$d.Content.Find.Execute("27/09/2024", $true, $false, $false, $false, $false,
                         $true, 1, $false, "22/10/2024", 2)

# Replace the student name placeholder "YUG" wherever it appears (3 occurrences)
$d.Content.Find.Execute("YUG", $true, $false, $false, $false, $false,
                         $true, 1, $false, "PARV SHAH (TESTING)", 2)

# Replace the enrolment number
$d.Content.Find.Execute("56757", $true, $false, $false, $false, $false,
                         $true, 1, $false, "226540307098", 2)

# Replace the admission quota/committee text
$d.Content.Find.Execute("Vacant Quota(Government)", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Central Admissions Committee", 2)
